# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    holding the per-fund holdings detail for the quarter.
# 2. Insert a new summary row at the top of the "总计" sheet's data for
#    "2022-Q1", shifting the existing rows down and renumbering the index
#    column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")

$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Re-resolve "总计" by name AFTER the insert: worksheet references in this
# host track sheet position, and inserting a sheet shifts "总计" from
# index 4 to index 5, so a handle captured beforehand would silently end
# up aliasing the new "2022-Q1" sheet instead.
$total = $wb.Worksheets.Item("总计")

# Borrow the existing bold/centered/bordered header style (used by the
# other quarter sheets) via copy/paste-special so we reuse style index 2
# instead of synthesizing a new one. Copy from "2021-Q4" (not "总计")
# because it actually has populated cells all the way out to column H --
# pasting from a narrower populated range leaves the extra columns
# unstyled.
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$q4.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0

$q4.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1

# Row 2: 南方绩优成长混合A (202003)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "202003"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "南方绩优成长混合A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.96"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "65.45"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2.13"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1.0428"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 9

# Row 3: 南方绩优成长混合C (006540)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "006540"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "南方绩优成长混合C"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.38"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "65.45"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2.13"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.0081"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" sheet
# ---------------------------------------------------------------------
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.05

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
